$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("51÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "53÷5=", 1) | Out-Null
$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("86÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "82÷5=", 1) | Out-Null
$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("52÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "80÷7=", 1) | Out-Null
$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("53÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "19÷5=", 1) | Out-Null
$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("81÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "50÷5=", 1) | Out-Null

$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("51÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "71÷7=", 1) | Out-Null
$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("72÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "80÷9=", 1) | Out-Null
$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("45÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "86÷9=", 1) | Out-Null
$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("57÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "79÷4=", 1) | Out-Null
$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("33÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "76÷7=", 1) | Out-Null

$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("38÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "29÷6=", 1) | Out-Null
$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("82÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "57÷3=", 1) | Out-Null
$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("26÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "61÷8=", 1) | Out-Null
$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("41÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "16÷6=", 1) | Out-Null
$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("30÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "17÷2=", 1) | Out-Null

$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("49÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "44÷9=", 1) | Out-Null
$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("94÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "11÷4=", 1) | Out-Null
$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("51÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "79÷6=", 1) | Out-Null
$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("76÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "41÷3=", 1) | Out-Null
$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("93÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "46÷9=", 1) | Out-Null

$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("39÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "18÷2=", 1) | Out-Null
$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("48÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "57÷7=", 1) | Out-Null
$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("86÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "23÷7=", 1) | Out-Null
$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("47÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "90÷2=", 1) | Out-Null
$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("65÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "86÷7=", 1) | Out-Null

Write-Host "Replacements complete"
